# Refresh the "cryptos" table (Price / Volume(1h) columns, plus two
# coins that swapped rank and had their Coin/Link/Price/Volume cells
# re-written in place) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.241.05"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").Value = "2.520.24"
$ws.Range("E3").Value = "  +1.38%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "590.63"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").Value = "177.85"
$ws.Range("E6").Value = "  +3.77%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("D9").Value = "0.145"
$ws.Range("E9").Value = "  +5.38%  "

$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("D11").Value = "0.342"
$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("D12").Value = "4.96"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").Value = "25.87"
$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").Value = "67.996.43"
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.850.11"
$ws.Range("E16").Value = "  -2.13%  "

$ws.Range("D17").Value = "2.429.25"
$ws.Range("E17").Value = "  -1.44%  "

$ws.Range("D18").Value = "11.09"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  +2.17%  "

$ws.Range("D20").Value = "353.58"
$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("D21").Value = "4.13"
$ws.Range("E21").Value = "  +2.45%  "

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").Value = "70.98"
$ws.Range("E23").Value = "  +3.74%  "

$ws.Range("D24").Value = "4.35"
$ws.Range("E24").Value = "  +2.98%  "

$ws.Range("D25").Value = "1.78"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").Value = "9.22"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").Value = "2.643.21"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "0.0₃0925"
$ws.Range("E29").Value = "  +2.10%  "

$ws.Range("D30").Value = "511.39"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("E32").Value = "  +3.74%  "

$ws.Range("E33").Value = "  +1.34%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  +3.91%  "

$ws.Range("D36").Value = "165.08"
$ws.Range("E36").Value = "  +2.67%  "

$ws.Range("D37").Value = "18.49"
$ws.Range("E37").Value = "  +1.26%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("E39").Value = "  +0.99%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").Value = "  +3.82%  "

$ws.Range("E42").Value = "  +2.15%  "

$ws.Range("D43").Value = "0.331"
$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("D44").Value = "2.51"
$ws.Range("E44").Value = "  +5.91%  "

$ws.Range("D45").Value = "147.78"
$ws.Range("E45").Value = "  +3.37%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0262"
$ws.Range("E47").Value = "  +4.17%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "0.524"
$ws.Range("E48").Value = "  +1.78%  "

$ws.Range("D49").Value = "0.0746"
$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("E50").Value = "  +2.58%  "

# D51's new price "0.590" has a trailing zero; assigning it to the cell
# directly would auto-detect it as a number and silently normalize it to
# "0.59". Write it into a scratch cell pre-formatted as Text, then
# paste-values-only onto D51 so the literal text is kept without leaving
# D51's own (default) number format/style touched.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "0.590"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E51").Value = "  +1.27%  "
